$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.67568533333333
$ws.Range("H2").Value = 83.027056
$ws.Range("I2").Value = 0.151580065893459
$ws.Range("J2").Value = 0.151580065893459
$ws.Range("M2").Value = 1.424886333333333
$ws.Range("N2").Value = 4.274659
$ws.Range("O2").Value = 0.3775790263072122
$ws.Range("P2").Value = 0.3775790263072122
$ws.Range("Q2").Value = 39.43470579710044
$ws.Range("R2").Value = 354.912352173904
$ws.Range("S2").Value = 0.0572334536876353
$ws.Range("T2").Value = 0.05723345368763531
$ws.Range("G3").Value = 27.67568533333333
$ws.Range("H3").Value = 83.027056
$ws.Range("I3").Value = 0.151580065893459
$ws.Range("J3").Value = 0.151580065893459
$ws.Range("O3").Value = 0.09562220712767076
$ws.Range("P3").Value = 0.09562220712767076
$ws.Range("Q3").Value = 9.986872530046224
$ws.Range("R3").Value = 89.88185277041602
$ws.Range("S3").Value = 0.01449442045729032
$ws.Range("T3").Value = 0.01449442045729032
$ws.Range("G4").Value = 27.67568533333333
$ws.Range("H4").Value = 83.027056
$ws.Range("I4").Value = 0.151580065893459
$ws.Range("J4").Value = 0.151580065893459
$ws.Range("M4").Value = 0.6353876666666666
$ws.Range("N4").Value = 1.906163
$ws.Range("O4").Value = 0.1683706629050024
$ws.Range("P4").Value = 0.1683706629050024
$ws.Range("Q4").Value = 17.58478912734755
$ws.Range("R4").Value = 158.263102146128
$ws.Range("S4").Value = 0.02552163617766562
$ws.Range("T4").Value = 0.02552163617766563
$ws.Range("G5").Value = 27.67568533333333
$ws.Range("H5").Value = 83.027056
$ws.Range("I5").Value = 0.151580065893459
$ws.Range("J5").Value = 0.151580065893459
$ws.Range("M5").Value = 0.4155976666666667
$ws.Range("N5").Value = 1.246793
$ws.Range("O5").Value = 0.110128758094306
$ws.Range("P5").Value = 0.110128758094306
$ws.Range("Q5").Value = 11.50195024793422
$ws.Range("R5").Value = 103.517552231408
$ws.Range("S5").Value = 0.01669332440869971
$ws.Range("T5").Value = 0.01669332440869971
$ws.Range("G6").Value = 27.67568533333333
$ws.Range("H6").Value = 83.027056
$ws.Range("I6").Value = 0.151580065893459
$ws.Range("J6").Value = 0.151580065893459
$ws.Range("M6").Value = 0.937018
$ws.Range("N6").Value = 2.811054
$ws.Range("O6").Value = 0.2482993455658087
$ws.Range("P6").Value = 0.2482993455658087
$ws.Range("Q6").Value = 25.93261531966933
$ws.Range("R6").Value = 233.393537877024
$ws.Range("S6").Value = 0.03763723116216802
$ws.Range("T6").Value = 0.03763723116216802
$ws.Range("I7").Value = 0.2439851776203359
$ws.Range("J7").Value = 0.243985177620336
$ws.Range("M7").Value = 1.424886333333333
$ws.Range("N7").Value = 4.274659
$ws.Range("O7").Value = 0.3775790263072122
$ws.Range("P7").Value = 0.3775790263072122
$ws.Range("Q7").Value = 63.47459767614755
$ws.Range("R7").Value = 571.2713790853279
$ws.Range("S7").Value = 0.09212368579927867
$ws.Range("T7").Value = 0.09212368579927867
$ws.Range("I8").Value = 0.2439851776203359
$ws.Range("J8").Value = 0.243985177620336
$ws.Range("O8").Value = 0.09562220712767076
$ws.Range("P8").Value = 0.09562220712767076
$ws.Range("S8").Value = 0.02333040119049331
$ws.Range("T8").Value = 0.02333040119049331
$ws.Range("I9").Value = 0.2439851776203359
$ws.Range("J9").Value = 0.243985177620336
$ws.Range("M9").Value = 0.6353876666666666
$ws.Range("N9").Value = 1.906163
$ws.Range("O9").Value = 0.1683706629050024
$ws.Range("P9").Value = 0.1683706629050024
$ws.Range("Q9").Value = 28.30469741098844
$ws.Range("R9").Value = 254.742276698896
$ws.Range("S9").Value = 0.04107994609493071
$ws.Range("T9").Value = 0.04107994609493072
$ws.Range("I10").Value = 0.2439851776203359
$ws.Range("J10").Value = 0.243985177620336
$ws.Range("M10").Value = 0.4155976666666667
$ws.Range("N10").Value = 1.246793
$ws.Range("O10").Value = 0.110128758094306
$ws.Range("P10").Value = 0.110128758094306
$ws.Range("Q10").Value = 18.51368356176178
$ws.Range("R10").Value = 166.623152055856
$ws.Range("S10").Value = 0.02686978460474626
$ws.Range("T10").Value = 0.02686978460474627
$ws.Range("I11").Value = 0.2439851776203359
$ws.Range("J11").Value = 0.243985177620336
$ws.Range("M11").Value = 0.937018
$ws.Range("N11").Value = 2.811054
$ws.Range("O11").Value = 0.2482993455658087
$ws.Range("P11").Value = 0.2482993455658087
$ws.Range("Q11").Value = 41.74146328301867
$ws.Range("R11").Value = 375.673169547168
$ws.Range("S11").Value = 0.060581359930887
$ws.Range("T11").Value = 0.06058135993088701
$ws.Range("G12").Value = 54.059897
$ws.Range("H12").Value = 162.179691
$ws.Range("I12").Value = 0.2960867147735651
$ws.Range("J12").Value = 0.2960867147735651
$ws.Range("M12").Value = 1.424886333333333
$ws.Range("N12").Value = 4.274659
$ws.Range("O12").Value = 0.3775790263072122
$ws.Range("P12").Value = 0.3775790263072122
$ws.Range("Q12").Value = 77.02920841670766
$ws.Range("R12").Value = 693.2628757503689
$ws.Range("S12").Value = 0.111796133466704
$ws.Range("T12").Value = 0.111796133466704
$ws.Range("G13").Value = 54.059897
$ws.Range("H13").Value = 162.179691
$ws.Range("I13").Value = 0.2960867147735651
$ws.Range("J13").Value = 0.2960867147735651
$ws.Range("O13").Value = 0.09562220712767076
$ws.Range("P13").Value = 0.09562220712767076
$ws.Range("Q13").Value = 19.50771205207234
$ws.Range("R13").Value = 175.569408468651
$ws.Range("S13").Value = 0.02831246516782942
$ws.Range("T13").Value = 0.02831246516782942
$ws.Range("G14").Value = 54.059897
$ws.Range("H14").Value = 162.179691
$ws.Range("I14").Value = 0.2960867147735651
$ws.Range("J14").Value = 0.2960867147735651
$ws.Range("M14").Value = 0.6353876666666666
$ws.Range("N14").Value = 1.906163
$ws.Range("O14").Value = 0.1683706629050024
$ws.Range("P14").Value = 0.1683706629050024
$ws.Range("Q14").Value = 34.34899181507033
$ws.Range("R14").Value = 309.140926335633
$ws.Range("S14").Value = 0.04985231644378951
$ws.Range("T14").Value = 0.04985231644378951
$ws.Range("G15").Value = 54.059897
$ws.Range("H15").Value = 162.179691
$ws.Range("I15").Value = 0.2960867147735651
$ws.Range("J15").Value = 0.2960867147735651
$ws.Range("M15").Value = 0.4155976666666667
$ws.Range("N15").Value = 1.246793
$ws.Range("O15").Value = 0.110128758094306
$ws.Range("P15").Value = 0.110128758094306
$ws.Range("Q15").Value = 22.46716705344033
$ws.Range("R15").Value = 202.204503480963
$ws.Range("S15").Value = 0.03260766218623573
$ws.Range("T15").Value = 0.03260766218623573
$ws.Range("G16").Value = 54.059897
$ws.Range("H16").Value = 162.179691
$ws.Range("I16").Value = 0.2960867147735651
$ws.Range("J16").Value = 0.2960867147735651
$ws.Range("M16").Value = 0.937018
$ws.Range("N16").Value = 2.811054
$ws.Range("O16").Value = 0.2482993455658087
$ws.Range("P16").Value = 0.2482993455658087
$ws.Range("Q16").Value = 50.655096567146
$ws.Range("R16").Value = 455.895869104314
$ws.Range("S16").Value = 0.07351813750900646
$ws.Range("T16").Value = 0.07351813750900646
$ws.Range("G17").Value = 11.41370466666667
$ws.Range("H17").Value = 34.241114
$ws.Range("I17").Value = 0.0625129995743248
$ws.Range("J17").Value = 0.0625129995743248
$ws.Range("M17").Value = 1.424886333333333
$ws.Range("N17").Value = 4.274659
$ws.Range("O17").Value = 0.3775790263072122
$ws.Range("P17").Value = 0.3775790263072122
$ws.Range("Q17").Value = 16.26323179223622
$ws.Range("R17").Value = 146.369086130126
$ws.Range("S17").Value = 0.02360359751081673
$ws.Range("T17").Value = 0.02360359751081673
$ws.Range("G18").Value = 11.41370466666667
$ws.Range("H18").Value = 34.241114
$ws.Range("I18").Value = 0.0625129995743248
$ws.Range("J18").Value = 0.0625129995743248
$ws.Range("O18").Value = 0.09562220712767076
$ws.Range("P18").Value = 0.09562220712767076
$ws.Range("Q18").Value = 4.118677179217111
$ws.Range("R18").Value = 37.06809461295401
$ws.Range("S18").Value = 0.00597763099346808
$ws.Range("T18").Value = 0.00597763099346808
$ws.Range("G19").Value = 11.41370466666667
$ws.Range("H19").Value = 34.241114
$ws.Range("I19").Value = 0.0625129995743248
$ws.Range("J19").Value = 0.0625129995743248
$ws.Range("M19").Value = 0.6353876666666666
$ws.Range("N19").Value = 1.906163
$ws.Range("O19").Value = 0.1683706629050024
$ws.Range("P19").Value = 0.1683706629050024
$ws.Range("Q19").Value = 7.252127176175778
$ws.Range("R19").Value = 65.269144585582
$ws.Range("S19").Value = 0.0105253551785092
$ws.Range("T19").Value = 0.0105253551785092
$ws.Range("G20").Value = 11.41370466666667
$ws.Range("H20").Value = 34.241114
$ws.Range("I20").Value = 0.0625129995743248
$ws.Range("J20").Value = 0.0625129995743248
$ws.Range("M20").Value = 0.4155976666666667
$ws.Range("N20").Value = 1.246793
$ws.Range("O20").Value = 0.110128758094306
$ws.Range("P20").Value = 0.110128758094306
$ws.Range("Q20").Value = 4.743509027489112
$ws.Range("R20").Value = 42.69158124740201
$ws.Range("S20").Value = 0.00688447900787027
$ws.Range("T20").Value = 0.00688447900787027
$ws.Range("G21").Value = 11.41370466666667
$ws.Range("H21").Value = 34.241114
$ws.Range("I21").Value = 0.0625129995743248
$ws.Range("J21").Value = 0.0625129995743248
$ws.Range("M21").Value = 0.937018
$ws.Range("N21").Value = 2.811054
$ws.Range("O21").Value = 0.2482993455658087
$ws.Range("P21").Value = 0.2482993455658087
$ws.Range("Q21").Value = 10.69484671935067
$ws.Range("R21").Value = 96.25362047415601
$ws.Range("S21").Value = 0.01552193688366052
$ws.Range("T21").Value = 0.01552193688366052
$ws.Range("G22").Value = 44.88488133333333
$ws.Range("H22").Value = 134.654644
$ws.Range("I22").Value = 0.2458350421383152
$ws.Range("J22").Value = 0.2458350421383153
$ws.Range("M22").Value = 1.424886333333333
$ws.Range("N22").Value = 4.274659
$ws.Range("O22").Value = 0.3775790263072122
$ws.Range("P22").Value = 0.3775790263072122
$ws.Range("Q22").Value = 63.95585398515511
$ws.Range("R22").Value = 575.602685866396
$ws.Range("S22").Value = 0.09282215584277755
$ws.Range("T22").Value = 0.09282215584277755
$ws.Range("G23").Value = 44.88488133333333
$ws.Range("H23").Value = 134.654644
$ws.Range("I23").Value = 0.2458350421383152
$ws.Range("J23").Value = 0.2458350421383153
$ws.Range("O23").Value = 0.09562220712767076
$ws.Range("P23").Value = 0.09562220712767076
$ws.Range("Q23").Value = 16.19687400703156
$ws.Range("R23").Value = 145.771866063284
$ws.Range("S23").Value = 0.02350728931858965
$ws.Range("T23").Value = 0.02350728931858965
$ws.Range("G24").Value = 44.88488133333333
$ws.Range("H24").Value = 134.654644
$ws.Range("I24").Value = 0.2458350421383152
$ws.Range("J24").Value = 0.2458350421383153
$ws.Range("M24").Value = 0.6353876666666666
$ws.Range("N24").Value = 1.906163
$ws.Range("O24").Value = 0.1683706629050024
$ws.Range("P24").Value = 0.1683706629050024
$ws.Range("Q24").Value = 28.51930001899688
$ws.Range("R24").Value = 256.673700170972
$ws.Range("S24").Value = 0.04139140901010732
$ws.Range("T24").Value = 0.04139140901010733
$ws.Range("G25").Value = 44.88488133333333
$ws.Range("H25").Value = 134.654644
$ws.Range("I25").Value = 0.2458350421383152
$ws.Range("J25").Value = 0.2458350421383153
$ws.Range("M25").Value = 0.4155976666666667
$ws.Range("N25").Value = 1.246793
$ws.Range("O25").Value = 0.110128758094306
$ws.Range("P25").Value = 0.110128758094306
$ws.Range("Q25").Value = 18.65405195074356
$ws.Range("R25").Value = 167.886467556692
$ws.Range("S25").Value = 0.02707350788675404
$ws.Range("T25").Value = 0.02707350788675405
$ws.Range("G26").Value = 44.88488133333333
$ws.Range("H26").Value = 134.654644
$ws.Range("I26").Value = 0.2458350421383152
$ws.Range("J26").Value = 0.2458350421383153
$ws.Range("M26").Value = 0.937018
$ws.Range("N26").Value = 2.811054
$ws.Range("O26").Value = 0.2482993455658087
$ws.Range("P26").Value = 0.2482993455658087
$ws.Range("Q26").Value = 42.05794173719733
$ws.Range("R26").Value = 378.521475634776
$ws.Range("S26").Value = 0.06104068008008667
$ws.Range("T26").Value = 0.06104068008008667
